# "change the login in datacenter frame"
# Split the single LoginDetails sheet into:
#  - LoginDetails        : now holds only invalid-credential test rows + an ErrorMessage column
#  - ValidLoginDetails    : a new sheet holding the one valid Admin/Qedge123!@# credential pair

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new worksheet right after LoginDetails --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ValidLoginDetails"

# --- Update LoginDetails (sheet1) ------------------------------------------
# Add the ErrorMessage column. Values are assigned in this specific order so
# new shared strings land at the expected indices (7,8,9,10).
$ws1.Range("C2").Value = "Password cannot be empty"
$ws1.Range("C1").Value = "ErrorMessage"
$ws1.Range("C3").Value = "Username cannot be empty"
$ws1.Range("C4").Value = "Username cannot be empty"
$ws1.Range("C5").Value = "Invalid credentials"
$ws1.Range("C6").Value = "Invalid credentials"

# Row 6 now holds what used to be row 7 (admin/password); the old valid
# Admin/Qedge123!@# pair (old row 6) has moved to the new sheet.
$ws1.Range("A6").Value = "admin"
$ws1.Range("B6").Value = "password"

# Old row 7 no longer exists (sheet now only spans rows 1-6).
$ws1.Range("A7:B7").ClearContents()

# --- Populate ValidLoginDetails (sheet2) -----------------------------------
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "Qedge123!@#"

# --- Selections / active cells to match the target worksheet views --------
$ws2.Range("F22").Select()
$ws1.Range("C1").Select()
